# "Generate Report for Handback"
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de translations have been handed back (in sync with en-US).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6003d8f39c71b90e3170c99c6045cfdf5371fc6/e2e/cf419707-8317-4710-b7fc-d771b2b6f4e3.md"
$mdName = "cf419707-8317-4710-b7fc-d771b2b6f4e3.md"

$statusText = "Handed back: in sync with en-US"

# --- Status text everywhere it is referenced -----------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsZhCn.Range("C2").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText

# --- zh-cn / de-de sheets: fill in target/handback info for row 2 ---------
# "Latest Target File" (I2) now points at the source markdown file.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdName)

# Apply the HyperLink look consistently (grouped by operation so the
# engine reuses a single deduplicated style/font across both sheets).
$wsZhCn.Range("I2").Style = "HyperLink"
$wsDeDe.Range("I2").Style = "HyperLink"
$wsZhCn.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsDeDe.Range("I2").Font.Color = 15570276

# "Latest Handback File" (J2) now points at the generated xliff file.
$wsZhCn.Range("J2").Value = "cf419707-8317-4710-b7fc-d771b2b6f4e3.8a31383afab559034b7f305dabeb116c632ef884.zh-cn.xlf"
$wsDeDe.Range("J2").Value = "cf419707-8317-4710-b7fc-d771b2b6f4e3.8a31383afab559034b7f305dabeb116c632ef884.de-de.xlf"

# "Latest Handback DateTime" (K2).
$wsZhCn.Range("K2").Value = "2016-09-06 15:41:59"
$wsDeDe.Range("K2").Value = "2016-09-06 15:42:34"

# --- Column width adjustments (content grew, columns widen) ---------------
# Helper values: Excel's ColumnWidth setter snaps to a 1/6-character grid
# (stored = round(ColumnWidth*6)/6 + 5/6), so we pick inputs that land on
# the closest achievable grid point to the target stored widths.
$wideWidth = 29.166666666666668   # -> stored width 30   (was ~17.22)
$fullWidth = 39.166666666666664   # -> stored width 40   (was ~18.65 / 21.71)

$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $fullWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $fullWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $fullWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $fullWidth
